$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.80267866345620553
$ws.Range("AN2").Value = 0.68344608205302371
$ws.Range("BH2").Value = 0.70364056235672012
$ws.Range("E3").Value = 0.87897915889162093
$ws.Range("B4").Value = 0.81831673521231885
$ws.Range("C4").Value = 0.99716132322386486
$ws.Range("F4").Value = 0.70298238169696703
$ws.Range("D5").Value = 0.89812434577910549
$ws.Range("AJ5").Value = 0.65267588491708928
$ws.Range("G6").Value = 0.67400823238869223
$ws.Range("H6").Value = 0.96744970575757039
$ws.Range("E7").Value = 0.99181165402860683
$ws.Range("I7").Value = 0.5994133865381841
$ws.Range("G8").Value = 0.98254588738927828
$ws.Range("J8").Value = 0.71559282962129012
$ws.Range("AV8").Value = 0.91138186866757054
$ws.Range("J9").Value = 0.90020443262553229
$ws.Range("S10").Value = 0.88468448064980687
$ws.Range("AZ10").Value = 0.89770312389180418
$ws.Range("I11").Value = 0.68589787438046224
$ws.Range("K12").Value = 0.61310733948352658
$ws.Range("AA12").Value = 0.88199191872085403
$ws.Range("BE12").Value = 0.95649154605226006
$ws.Range("A13").Value = 0.64305736457437934
$ws.Range("K13").Value = 0.90121002702044795
$ws.Range("L13").Value = 0.90061730832820142
$ws.Range("N13").Value = 0.89983998581049907
$ws.Range("P14").Value = 0.72079304398852173
$ws.Range("BI15").Value = 0.96549222255817391
$ws.Range("I16").Value = 0.71559405493855999
$ws.Range("P18").Value = 0.89410976290748589
$ws.Range("AL18").Value = 0.96020093571051368
$ws.Range("N19").Value = 0.69509697441083385
$ws.Range("T19").Value = 0.86748956362763674
$ws.Range("R20").Value = 0.61013148336524159
$ws.Range("W21").Value = 0.54850595409677327
$ws.Range("U22").Value = 0.94137135071010136
$ws.Range("V23").Value = 0.66576040914501866
$ws.Range("Y23").Value = 0.63632379789377502
$ws.Range("BO23").Value = 0.95634630204890558
$ws.Range("V24").Value = 0.94950092679113229
$ws.Range("Y24").Value = 0.95229474307312534
$ws.Range("Z25").Value = 0.92437537664845282
$ws.Range("X26").Value = 0.7823685626598319
$ws.Range("H27").Value = 0.94926817356885695
$ws.Range("Y27").Value = 0.82538797846795053
$ws.Range("Z28").Value = 0.94834683954358789
$ws.Range("AD28").Value = 0.95523570648527523
$ws.Range("AB29").Value = 0.886784793129749
$ws.Range("AD29").Value = 0.77092026634486766
$ws.Range("V30").Value = 0.99809481666813293
$ws.Range("S31").Value = 0.94983393886767797
$ws.Range("AC31").Value = 0.92112598203689644
$ws.Range("AD31").Value = 0.80862032262875738
$ws.Range("AK31").Value = 0.89083507908378734
$ws.Range("AG32").Value = 0.98544165265566364
$ws.Range("AH32").Value = 0.85090087859213093
$ws.Range("AI33").Value = 0.95052478515474226
$ws.Range("F34").Value = 0.63449251980192845
$ws.Range("AG34").Value = 0.7504976833050101
$ws.Range("AI34").Value = 0.77882600182979922
$ws.Range("J36").Value = 0.89378874088249272
$ws.Range("AK36").Value = 0.58919166040030635
$ws.Range("P37").Value = 0.90094397760983203
$ws.Range("AI37").Value = 0.70038306430141306
$ws.Range("AL37").Value = 0.97721603042906757
$ws.Range("AC38").Value = 0.93013887360982617
$ws.Range("AJ38").Value = 0.8955018258843489
$ws.Range("AL39").Value = 0.95540356079205158
$ws.Range("AN39").Value = 0.79075795476098298
$ws.Range("BF40").Value = 0.87767807772337858
$ws.Range("Q41").Value = 0.93572029774622323
$ws.Range("AM41").Value = 0.64538677273300715
$ws.Range("AQ41").Value = 0.76040033433001875
$ws.Range("AO42").Value = 0.84865995351547863
$ws.Range("AQ42").Value = 0.86208156531400992
$ws.Range("AR42").Value = 0.94848075856213243
$ws.Range("J43").Value = 0.82178566794644747
$ws.Range("AR43").Value = 0.75086421678115101
$ws.Range("BH44").Value = 0.61903726433600348
$ws.Range("F45").Value = 0.78704722184215115
$ws.Range("AQ45").Value = 0.83491185696599712
$ws.Range("AT45").Value = 0.99585242454474543
$ws.Range("BD46").Value = 0.86408744357310141
$ws.Range("G47").Value = 0.98677447461799095
$ws.Range("AT48").Value = 0.95186905608229022
$ws.Range("F49").Value = 0.60339500149336911
$ws.Range("AU49").Value = 0.80238241544238831
$ws.Range("AV49").Value = 0.68332134247138288
$ws.Range("AY49").Value = 0.88760112497844301
$ws.Range("AP50").Value = 0.68511107314536457
$ws.Range("BA51").Value = 0.83107035271078433
$ws.Range("AI52").Value = 0.54653732222919493
$ws.Range("AX52").Value = 0.8785882682133006
$ws.Range("AY52").Value = 0.94607243890466242
$ws.Range("BB53").Value = 0.9506935022377061
$ws.Range("AZ54").Value = 0.8478812019194657
$ws.Range("BC54").Value = 0.9722561431900133
$ws.Range("BD54").Value = 0.94006060615819942
$ws.Range("O55").Value = 0.71016011452044325
$ws.Range("AA55").Value = 0.9337960335012685
$ws.Range("AP55").Value = 0.84872377682993738
$ws.Range("BA55").Value = 0.86221584419481478
$ws.Range("F56").Value = 0.83039007255685626
$ws.Range("BF56").Value = 0.90290120842776234
$ws.Range("AZ57").Value = 0.91375464965352315
$ws.Range("BK57").Value = 0.97497747192231521
$ws.Range("BE58").Value = 0.57806070876041371
$ws.Range("BG58").Value = 0.96647693386440103
$ws.Range("Z59").Value = 0.6937868098081591
$ws.Range("AR59").Value = 0.57195670452949621
$ws.Range("BH59").Value = 0.68805839938869251
$ws.Range("BI59").Value = 0.9340000022475704
$ws.Range("A60").Value = 0.86676405400017065
$ws.Range("BF60").Value = 0.93042246103404924
$ws.Range("BJ60").Value = 0.97047910111461577
$ws.Range("BH61").Value = 0.79357126750368701
$ws.Range("BJ61").Value = 0.72572534226836849
$ws.Range("BK61").Value = 0.96050072941946185
$ws.Range("T62").Value = 0.86171440723965342
$ws.Range("BL62").Value = 0.67562809776443755
$ws.Range("T63").Value = 0.83875536922228755
$ws.Range("BJ63").Value = 0.78565300823805484
$ws.Range("BN64").Value = 0.88205210060849371
$ws.Range("BK65").Value = 0.88530629483226653
$ws.Range("BL65").Value = 0.90580729722643127
$ws.Range("AA66").Value = 0.88196736336956505
$ws.Range("BD66").Value = 0.85927841057521803
$ws.Range("A67").Value = 0.89524451396117311
$ws.Range("Q67").Value = 0.88485321420400698
$ws.Range("AM67").Value = 0.74754742679900288
$ws.Range("BM67").Value = 0.92310715131313059
$ws.Range("BN67").Value = 0.87617303428597393
$ws.Range("B68").Value = 0.78448566156468269
$ws.Range("AE68").Value = 0.70243108661681797
$ws.Range("BM68").Value = 0.86077351763022691
